$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4 (shifts existing rows 4-29 down to 6-31)
$ws.Rows("4:5").Insert()

# Copy formatting (border/font/alignment) from column-A style onto the new rows
$ws.Range("A6").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 4: Holden
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Holden"
$ws.Cells.Item(4, 3).Value = 1.021647457291363
$ws.Cells.Item(4, 4).Value = 0.9282671978906109
$ws.Cells.Item(4, 5).Value = 1.011311760564904
$ws.Cells.Item(4, 6).Value = 0.9282671978906109
$ws.Cells.Item(4, 7).Value = 1.011311760564904
$ws.Cells.Item(4, 8).Value = 1.05987875461917
$ws.Cells.Item(4, 9).Value = 0.9580522655597753
$ws.Cells.Item(4, 10).Value = 1.01863413406255
$ws.Cells.Item(4, 11).Value = 1.011311760564904
$ws.Cells.Item(4, 12).Value = 1.021647457291363
$ws.Cells.Item(4, 13).Value = 0.9749573275909871
$ws.Cells.Item(4, 14).Value = 0.9749573275909871
$ws.Cells.Item(4, 15).Value = 0.9693223069139165
$ws.Cells.Item(4, 16).Value = 0.9870754719156259
$ws.Cells.Item(4, 17).Value = 0.9870754719156259
$ws.Cells.Item(4, 18).Value = 0.9931345440779453
$ws.Cells.Item(4, 19).Value = 0.9931345440779453
$ws.Cells.Item(4, 20).Value = 0.9996319283313956

# New row 5: Rizzie Spiral
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$ws.Cells.Item(5, 3).Value = 1.014507662776878
$ws.Cells.Item(5, 4).Value = 0.9496899277255813
$ws.Cells.Item(5, 5).Value = 1.009129976083197
$ws.Cells.Item(5, 6).Value = 0.9496899277255813
$ws.Cells.Item(5, 7).Value = 1.009129976083197
$ws.Cells.Item(5, 8).Value = 1.039787811668401
$ws.Cells.Item(5, 9).Value = 0.9710421488710952
$ws.Cells.Item(5, 10).Value = 1.012939825723358
$ws.Cells.Item(5, 11).Value = 1.009129976083197
$ws.Cells.Item(5, 12).Value = 1.014507662776878
$ws.Cells.Item(5, 13).Value = 0.9820987952512296
$ws.Cells.Item(5, 14).Value = 0.9820987952512296
$ws.Cells.Item(5, 15).Value = 0.9784132464578515
$ws.Cells.Item(5, 16).Value = 0.9911091888618854
$ws.Cells.Item(5, 17).Value = 0.9911091888618854
$ws.Cells.Item(5, 18).Value = 0.9956143856672133
$ws.Cells.Item(5, 19).Value = 0.9956143856672133
$ws.Cells.Item(5, 20).Value = 0.9995162254747516

# Fix the A-column sequence number for all shifted rows (6-31)
for ($r = 6; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Rename "Thomas Hex" -> "Matthies Hex" (now at row 11 after the shift)
$ws.Cells.Item(11, 2).Value = "Matthies Hex"

# Update dimension to reflect the two extra rows
$ws.UsedRange | Out-Null
